$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 1.030562134411285
$ws.Range("D2").Value = 1.033206538173445
$ws.Range("E2").Value = 1.038875712721936
$ws.Range("F2").Value = 1.046786514459898
$ws.Range("J2").Value = 1.035702694217099
$ws.Range("K2").Value = 1.036009554240408
$ws.Range("L2").Value = 1.041662500047411
$ws.Range("M2").Value = 1.049550974837245
$ws.Range("N2").Value = 1.037173510839646
$ws.Range("C3").Value = 1.03215381108919
$ws.Range("D3").Value = 1.03471179359483
$ws.Range("E3").Value = 1.04033598946218
$ws.Range("F3").Value = 1.048434650840371
$ws.Range("J3").Value = 1.036932603174502
$ws.Range("K3").Value = 1.0373218830727
$ws.Range("L3").Value = 1.042931164785066
$ws.Range("M3").Value = 1.051008653747325
$ws.Range("N3").Value = 1.038405166408842
$ws.Range("C4").Value = 1.033182472109806
$ws.Range("D4").Value = 1.035684889488324
$ws.Range("E4").Value = 1.041279939445193
$ws.Range("F4").Value = 1.049500404049907
$ws.Range("J4").Value = 1.037726867828673
$ws.Range("K4").Value = 1.038169651329701
$ws.Range("L4").Value = 1.043750624640977
$ws.Range("M4").Value = 1.051950700288226
$ws.Range("N4").Value = 1.039200559009926
$ws.Range("C5").Value = 1.033614629190545
$ws.Range("D5").Value = 1.036093770649095
$ws.Range("E5").Value = 1.041676557336047
$ws.Range("F5").Value = 1.049948288000722
$ws.Range("J5").Value = 1.038060408364389
$ws.Range("K5").Value = 1.038525726403061
$ws.Range("L5").Value = 1.044094785471044
$ws.Range("M5").Value = 1.052346465061968
$ws.Range("N5").Value = 1.03953457321146
$ws.Range("C6").Value = 1.033687173422782
$ws.Range("D6").Value = 1.036162411617872
$ws.Range("E6").Value = 1.041743138632082
$ws.Range("F6").Value = 1.05002348064651
$ws.Range("J6").Value = 1.038116389898298
$ws.Range("K6").Value = 1.03858549403571
$ws.Range("L6").Value = 1.044152551856734
$ws.Range("M6").Value = 1.052412900023368
$ws.Range("N6").Value = 1.039590634245568
$ws.Range("C7").Value = 1.033188247747824
$ws.Range("D7").Value = 1.035690353785499
$ws.Range("E7").Value = 1.041285239921016
$ws.Range("F7").Value = 1.049506389307972
$ws.Range("J7").Value = 1.037731326050825
$ws.Range("K7").Value = 1.038174410495421
$ws.Range("L7").Value = 1.043755224660985
$ws.Range("M7").Value = 1.051955989572085
$ws.Range("N7").Value = 1.039205023563264
$ws.Range("C8").Value = 1.031100315779474
$ws.Range("D8").Value = 1.033715438695404
$ws.Range("E8").Value = 1.039369420912191
$ws.Range("F8").Value = 1.047343660600443
$ws.Range("J8").Value = 1.036118676871632
$ws.Range("K8").Value = 1.036453355580747
$ws.Range("L8").Value = 1.042091556089259
$ws.Range("M8").Value = 1.050043851354095
$ws.Range("N8").Value = 1.037590084237263
$ws.Range("C9").Value = 1.027411026080917
$ws.Range("D9").Value = 1.030228081309539
$ws.Range("E9").Value = 1.03598589458151
$ws.Range("F9").Value = 1.043526873432085
$ws.Range("J9").Value = 1.033264644288842
$ws.Range("K9").Value = 1.033409603510073
$ws.Range("L9").Value = 1.039148518373584
$ws.Range("M9").Value = 1.046665091508326
$ws.Range("N9").Value = 1.034731998600843
$ws.Range("C10").Value = 1.024944108354034
$ws.Range("D10").Value = 1.027897741601417
$ws.Range("E10").Value = 1.033724593018178
$ws.Range("F10").Value = 1.040977882752943
$ws.Range("J10").Value = 1.031353204922554
$ws.Range("K10").Value = 1.031372540575154
$ws.Range("L10").Value = 1.037178341101313
$ws.Range("M10").Value = 1.044405789712116
$ws.Range("N10").Value = 1.03281784477141
$ws.Range("C11").Value = 1.023874016672767
$ws.Range("D11").Value = 1.026887272776326
$ws.Range("E11").Value = 1.032743979980415
$ws.Range("F11").Value = 1.039872953408709
$ws.Range("J11").Value = 1.030523353244349
$ws.Range("K11").Value = 1.030488490769379
$ws.Range("L11").Value = 1.036323196263315
$ws.Range("M11").Value = 1.043425758413411
$ws.Range("N11").Value = 1.031986814608637
$ws.Range("C12").Value = 1.023476239755991
$ws.Range("D12").Value = 1.026511716740801
$ws.Range("E12").Value = 1.032379508388751
$ws.Range("F12").Value = 1.039462342069675
$ws.Range("J12").Value = 1.03021477216249
$ws.Range("K12").Value = 1.030159807381617
$ws.Range("L12").Value = 1.036005241231162
$ws.Range("M12").Value = 1.043061459799599
$ws.Range("N12").Value = 1.031677795306247
$ws.Range("C13").Value = 1.023561577894621
$ws.Range("D13").Value = 1.026592285036533
$ws.Range("E13").Value = 1.032457699267963
$ws.Range("F13").Value = 1.039550428436437
$ws.Range("J13").Value = 1.03028097925438
$ws.Range("K13").Value = 1.030230325167483
$ws.Range("L13").Value = 1.03607345811585
$ws.Range("M13").Value = 1.043139615531749
$ws.Range("N13").Value = 1.031744096419801
$ws.Range("C14").Value = 1.023841142440326
$ws.Range("D14").Value = 1.026856233803733
$ws.Range("E14").Value = 1.032713857357642
$ws.Range("F14").Value = 1.039839016116397
$ws.Range("J14").Value = 1.030497852766341
$ws.Range("K14").Value = 1.030461328033341
$ws.Range("L14").Value = 1.036296920502646
$ws.Range("M14").Value = 1.043395650970781
$ws.Range("N14").Value = 1.031961277917025
$ws.Range("C15").Value = 1.024013351645639
$ws.Range("D15").Value = 1.027018831506691
$ws.Range("E15").Value = 1.032871654273831
$ws.Range("F15").Value = 1.040016798814951
$ws.Range("J15").Value = 1.030631430744719
$ws.Range("K15").Value = 1.030603615442074
$ws.Range("L15").Value = 1.036434560891634
$ws.Range("M15").Value = 1.043553366593536
$ws.Range("N15").Value = 1.032095045591454
$ws.Range("C16").Value = 1.025015085781351
$ws.Range("D16").Value = 1.02796477238556
$ws.Range("E16").Value = 1.033789641602626
$ws.Range("F16").Value = 1.041051187050347
$ws.Range("J16").Value = 1.031408232641063
$ws.Range("K16").Value = 1.031431169375307
$ws.Range("L16").Value = 1.037235050401181
$ws.Range("M16").Value = 1.04447079366613
$ws.Range("N16").Value = 1.032872950635592
$ws.Range("C17").Value = 1.025642930172393
$ws.Range("D17").Value = 1.028557749781809
$ws.Range("E17").Value = 1.034365074542139
$ws.Range("F17").Value = 1.041699702452391
$ws.Range("J17").Value = 1.031894908533827
$ws.Range("K17").Value = 1.031949733322234
$ws.Range("L17").Value = 1.037736622525418
$ws.Range("M17").Value = 1.04504579880025
$ws.Range("N17").Value = 1.033360317663945
$ws.Range("C18").Value = 1.026008958657824
$ws.Range("D18").Value = 1.028903487495892
$ws.Range("E18").Value = 1.034700575159188
$ws.Range("F18").Value = 1.042077855533128
$ws.Range("J18").Value = 1.032178568048327
$ws.Range("K18").Value = 1.032252011703005
$ws.Range("L18").Value = 1.038028984345274
$ws.Range("M18").Value = 1.045381022494943
$ws.Range("N18").Value = 1.033644380007492
$ws.Range("C19").Value = 1.026133734307175
$ws.Range("D19").Value = 1.029021352364181
$ws.Range("E19").Value = 1.034814948799851
$ws.Range("F19").Value = 1.042206776794718
$ws.Range("J19").Value = 1.032275253207596
$ws.Range("K19").Value = 1.032355048667241
$ws.Range("L19").Value = 1.038128639072815
$ws.Range("M19").Value = 1.04549529707875
$ws.Range("N19").Value = 1.033741202470777
$ws.Range("C20").Value = 1.025575587375707
$ws.Range("D20").Value = 1.028494143056313
$ws.Range("E20").Value = 1.03430335055181
$ws.Range("F20").Value = 1.041630134861107
$ws.Range("J20").Value = 1.031842714632578
$ws.Range("K20").Value = 1.031894116161544
$ws.Range("L20").Value = 1.03768282899494
$ws.Range("M20").Value = 1.044984123545516
$ws.Range("N20").Value = 1.03330804964137
$ws.Range("C21").Value = 1.023758825949248
$ws.Range("D21").Value = 1.026778513706073
$ws.Range("E21").Value = 1.032638431564003
$ws.Range("F21").Value = 1.039754039610723
$ws.Range("J21").Value = 1.030433998302629
$ws.Range("K21").Value = 1.030393312021792
$ws.Range("L21").Value = 1.036231125199006
$ws.Range("M21").Value = 1.043320262487956
$ws.Range("N21").Value = 1.031897332772651
$ws.Range("C22").Value = 1.022614830413416
$ws.Range("D22").Value = 1.025698534600439
$ws.Range("E22").Value = 1.031590305474058
$ws.Range("F22").Value = 1.038573351184445
$ws.Range("J22").Value = 1.029546326963721
$ws.Range("K22").Value = 1.029447910689423
$ws.Range("L22").Value = 1.03531654770955
$ws.Range("M22").Value = 1.042272551957321
$ws.Range("N22").Value = 1.031008400838624
$ws.Range("C23").Value = 1.023221451602037
$ws.Range("D23").Value = 1.026271178097185
$ws.Range("E23").Value = 1.032146065846889
$ws.Range("F23").Value = 1.039199365696136
$ws.Range("J23").Value = 1.030017086655981
$ws.Range("K23").Value = 1.02994925823704
$ws.Range("L23").Value = 1.035801559348376
$ws.Range("M23").Value = 1.042828115832002
$ws.Range("N23").Value = 1.031479829063643
$ws.Range("C24").Value = 1.025606017250856
$ws.Range("D24").Value = 1.02852288461791
$ws.Range("E24").Value = 1.034331241397905
$ws.Range("F24").Value = 1.041661569816516
$ws.Range("J24").Value = 1.031866299459806
$ws.Range("K24").Value = 1.03191924775334
$ws.Range("L24").Value = 1.037707136583794
$ws.Range("M24").Value = 1.045011992454706
$ws.Range("N24").Value = 1.033331667961758
$ws.Range("C25").Value = 1.02836605372284
$ws.Range("D25").Value = 1.031130566094369
$ws.Range("E25").Value = 1.036861570325309
$ws.Range("F25").Value = 1.044514348941243
$ws.Range("J25").Value = 1.034003988974581
$ws.Range("K25").Value = 1.034197842754593
$ws.Range("L25").Value = 1.03991076712819
$ws.Range("M25").Value = 1.047539741011015
$ws.Range("N25").Value = 1.03547239324084
